$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Йогурт"
$ws.Range("B4").Value = 110

$ws.Range("B5").Select()
